$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use row 283 as a style/structure template for the new rows (A283:M283 -> A284:M309)
$template = $ws.Range("A283:M283")
for ($r = 284; $r -le 309; $r++) {
    $template.Copy($ws.Range("A" + $r + ":M" + $r))
}

# Row 284 (NUMBER 283)
$ws.Cells.Item(284,1).Value = 283
$ws.Cells.Item(284,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(284,3).Value = "5:00 AM"
$ws.Cells.Item(284,4).Value = "E45159"
$ws.Cells.Item(284,5).Value = "Hurghada"
$ws.Cells.Item(284,6).Value = "(HRG)"
$ws.Cells.Item(284,7).Value = "Enter Air "
$ws.Cells.Item(284,8).Value = "B738"
$ws.Cells.Item(284,9).Value = "(SP-ESH)"
$ws.Cells.Item(284,10).Value = "5:01 AM"
$ws.Cells.Item(284,12).Value = "0 hours, 1 minutes"

# Row 285 (NUMBER 284)
$ws.Cells.Item(285,1).Value = 284
$ws.Cells.Item(285,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(285,3).Value = "5:10 AM"
$ws.Cells.Item(285,4).Value = "BO625"
$ws.Cells.Item(285,5).Value = "Madrid"
$ws.Cells.Item(285,6).Value = "(MAD)"
$ws.Cells.Item(285,7).Value = "Bluebird Nordic "
$ws.Cells.Item(285,8).Value = "B734"
$ws.Cells.Item(285,9).Value = "(TF-BBN)"
$ws.Cells.Item(285,10).Value = "5:04 AM"
$ws.Cells.Item(285,12).Value = "0 hours, -6 minutes"

# Row 286 (NUMBER 285)
$ws.Cells.Item(286,1).Value = 285
$ws.Cells.Item(286,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(286,3).Value = "5:30 AM"
$ws.Cells.Item(286,4).Value = "E44831"
$ws.Cells.Item(286,5).Value = "Hurghada"
$ws.Cells.Item(286,6).Value = "(HRG)"
$ws.Cells.Item(286,7).Value = "Enter Air "
$ws.Cells.Item(286,8).Value = "B738"
$ws.Cells.Item(286,9).Value = "(SP-ENP)"
$ws.Cells.Item(286,10).Value = "5:27 AM"
$ws.Cells.Item(286,12).Value = "0 hours, -3 minutes"

# Row 287 (NUMBER 286)
$ws.Cells.Item(287,1).Value = 286
$ws.Cells.Item(287,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(287,3).Value = "6:05 AM"
$ws.Cells.Item(287,4).Value = "FR6367"
$ws.Cells.Item(287,5).Value = "Catania"
$ws.Cells.Item(287,6).Value = "(CTA)"
$ws.Cells.Item(287,7).Value = "Ryanair "
$ws.Cells.Item(287,8).Value = "B738"
$ws.Cells.Item(287,9).Value = "(SP-RSH)"
$ws.Cells.Item(287,10).Value = "6:09 AM"
$ws.Cells.Item(287,12).Value = "0 hours, 4 minutes"

# Row 288 (NUMBER 287)
$ws.Cells.Item(288,1).Value = 287
$ws.Cells.Item(288,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(288,3).Value = "6:15 AM"
$ws.Cells.Item(288,4).Value = "W61001"
$ws.Cells.Item(288,5).Value = "London"
$ws.Cells.Item(288,6).Value = "(LTN)"
$ws.Cells.Item(288,7).Value = "Wizz Air "
$ws.Cells.Item(288,8).Value = "A21N"
$ws.Cells.Item(288,9).Value = "(HA-LZF)"
$ws.Cells.Item(288,10).Value = "6:23 AM"
$ws.Cells.Item(288,12).Value = "0 hours, 8 minutes"

# Row 289 (NUMBER 288)
$ws.Cells.Item(289,1).Value = 288
$ws.Cells.Item(289,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(289,3).Value = "6:20 AM"
$ws.Cells.Item(289,4).Value = "W61019"
$ws.Cells.Item(289,5).Value = "Bristol"
$ws.Cells.Item(289,6).Value = "(BRS)"
$ws.Cells.Item(289,7).Value = "Wizz Air "
$ws.Cells.Item(289,8).Value = "A21N"
$ws.Cells.Item(289,9).Value = "(HA-LZJ)"
$ws.Cells.Item(289,10).Value = "6:25 AM"
$ws.Cells.Item(289,12).Value = "0 hours, 5 minutes"

# Row 290 (NUMBER 289)
$ws.Cells.Item(290,1).Value = 289
$ws.Cells.Item(290,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(290,3).Value = "6:40 AM"
$ws.Cells.Item(290,4).Value = "W61241"
$ws.Cells.Item(290,5).Value = "Athens"
$ws.Cells.Item(290,6).Value = "(ATH)"
$ws.Cells.Item(290,7).Value = "Wizz Air "
$ws.Cells.Item(290,8).Value = "A21N"
$ws.Cells.Item(290,9).Value = "(HA-LVT)"
$ws.Cells.Item(290,10).Value = "6:41 AM"
$ws.Cells.Item(290,12).Value = "0 hours, 1 minutes"

# Row 291 (NUMBER 290)
$ws.Cells.Item(291,1).Value = 290
$ws.Cells.Item(291,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(291,3).Value = "7:00 AM"
$ws.Cells.Item(291,4).Value = "3Z7110"
$ws.Cells.Item(291,5).Value = "Hurghada"
$ws.Cells.Item(291,6).Value = "(HRG)"
$ws.Cells.Item(291,7).Value = "Smartwings "
$ws.Cells.Item(291,8).Value = "B738"
$ws.Cells.Item(291,9).Value = "(OK-TSF)"
$ws.Cells.Item(291,10).Value = "6:16 AM"
$ws.Cells.Item(291,12).Value = "0 hours, -44 minutes"

# Row 292 (NUMBER 291)
$ws.Cells.Item(292,1).Value = 291
$ws.Cells.Item(292,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(292,3).Value = "7:25 AM"
$ws.Cells.Item(292,4).Value = "W61071"
$ws.Cells.Item(292,5).Value = "Eindhoven"
$ws.Cells.Item(292,6).Value = "(EIN)"
$ws.Cells.Item(292,7).Value = "Wizz Air "
$ws.Cells.Item(292,8).Value = "A321"
$ws.Cells.Item(292,9).Value = "(HA-LTC)"
$ws.Cells.Item(292,10).Value = "7:43 AM"
$ws.Cells.Item(292,12).Value = "0 hours, 18 minutes"

# Row 293 (NUMBER 292)
$ws.Cells.Item(293,1).Value = 292
$ws.Cells.Item(293,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(293,3).Value = "8:00 AM"
$ws.Cells.Item(293,4).Value = "FR6892"
$ws.Cells.Item(293,5).Value = "Dortmund"
$ws.Cells.Item(293,6).Value = "(DTM)"
$ws.Cells.Item(293,7).Value = "Ryanair "
$ws.Cells.Item(293,8).Value = "B738"
$ws.Cells.Item(293,9).Value = "(SP-RSB)"
$ws.Cells.Item(293,10).Value = "8:00 AM"
$ws.Cells.Item(293,12).Value = "0 hours, 0 minutes"

# Row 294 (NUMBER 293)
$ws.Cells.Item(294,1).Value = 293
$ws.Cells.Item(294,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(294,3).Value = "8:50 AM"
$ws.Cells.Item(294,4).Value = "FR3409"
$ws.Cells.Item(294,5).Value = "Milan"
$ws.Cells.Item(294,6).Value = "(BGY)"
$ws.Cells.Item(294,7).Value = "Ryanair "
$ws.Cells.Item(294,8).Value = "B738"
$ws.Cells.Item(294,9).Value = "(9H-QAK)"
$ws.Cells.Item(294,10).Value = "8:53 AM"
$ws.Cells.Item(294,12).Value = "0 hours, 3 minutes"

# Row 295 (NUMBER 294)
$ws.Cells.Item(295,1).Value = 294
$ws.Cells.Item(295,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(295,3).Value = "9:35 AM"
$ws.Cells.Item(295,4).Value = "BO951"
$ws.Cells.Item(295,5).Value = "Paris"
$ws.Cells.Item(295,6).Value = "(CDG)"
$ws.Cells.Item(295,7).Value = "Bluebird Nordic "
$ws.Cells.Item(295,8).Value = "B734"
$ws.Cells.Item(295,9).Value = "(TF-BBJ)"
$ws.Cells.Item(295,10).Value = "9:38 AM"
$ws.Cells.Item(295,12).Value = "0 hours, 3 minutes"

# Row 296 (NUMBER 295)
$ws.Cells.Item(296,1).Value = 295
$ws.Cells.Item(296,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(296,3).Value = "9:50 AM"
$ws.Cells.Item(296,4).Value = "FR6391"
$ws.Cells.Item(296,5).Value = "London"
$ws.Cells.Item(296,6).Value = "(STN)"
$ws.Cells.Item(296,7).Value = "Ryanair "
$ws.Cells.Item(296,8).Value = "B38M"
$ws.Cells.Item(296,9).Value = "(EI-HMS)"
$ws.Cells.Item(296,10).Value = "10:05 AM"
$ws.Cells.Item(296,12).Value = "0 hours, 15 minutes"

# Row 297 (NUMBER 296)
$ws.Cells.Item(297,1).Value = 296
$ws.Cells.Item(297,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(297,3).Value = "11:00 AM"
$ws.Cells.Item(297,4).Value = "LO6543"
$ws.Cells.Item(297,5).Value = "Cancun"
$ws.Cells.Item(297,6).Value = "(CUN)"
$ws.Cells.Item(297,7).Value = "LOT "
$ws.Cells.Item(297,8).Value = "B788"
$ws.Cells.Item(297,9).Value = "(SP-LRC)"
$ws.Cells.Item(297,10).Value = "11:06 AM"
$ws.Cells.Item(297,12).Value = "0 hours, 6 minutes"

# Row 298 (NUMBER 297)
$ws.Cells.Item(298,1).Value = 297
$ws.Cells.Item(298,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(298,3).Value = "11:35 AM"
$ws.Cells.Item(298,4).Value = "W61093"
$ws.Cells.Item(298,5).Value = "Dortmund"
$ws.Cells.Item(298,6).Value = "(DTM)"
$ws.Cells.Item(298,7).Value = "Wizz Air "
$ws.Cells.Item(298,8).Value = "A321"
$ws.Cells.Item(298,9).Value = "(HA-LXN)"
$ws.Cells.Item(298,10).Value = "11:42 AM"
$ws.Cells.Item(298,12).Value = "0 hours, 7 minutes"

# Row 299 (NUMBER 298)
$ws.Cells.Item(299,1).Value = 298
$ws.Cells.Item(299,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(299,3).Value = "11:40 AM"
$ws.Cells.Item(299,4).Value = "LO3882"
$ws.Cells.Item(299,5).Value = "Warsaw"
$ws.Cells.Item(299,6).Value = "(WAW)"
$ws.Cells.Item(299,7).Value = "LOT "
$ws.Cells.Item(299,8).Value = "E170"
$ws.Cells.Item(299,9).Value = "(SP-LDI)"
$ws.Cells.Item(299,10).Value = "11:45 AM"
$ws.Cells.Item(299,12).Value = "0 hours, 5 minutes"

# Row 300 (NUMBER 299)
$ws.Cells.Item(300,1).Value = 299
$ws.Cells.Item(300,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(300,3).Value = "12:00 PM"
$ws.Cells.Item(300,4).Value = "UNKNOWN"
$ws.Cells.Item(300,5).Value = "Dammam"
$ws.Cells.Item(300,6).Value = "(DMM)"
$ws.Cells.Item(300,7).Value = "Enter Air "
$ws.Cells.Item(300,8).Value = "B738"
$ws.Cells.Item(300,9).Value = "(SP-ESC)"
$ws.Cells.Item(300,10).Value = "12:51 PM"
$ws.Cells.Item(300,12).Value = "0 hours, 51 minutes"

# Row 301 (NUMBER 300)
$ws.Cells.Item(301,1).Value = 300
$ws.Cells.Item(301,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(301,3).Value = "12:05 PM"
$ws.Cells.Item(301,4).Value = "FR6385"
$ws.Cells.Item(301,5).Value = "Athens"
$ws.Cells.Item(301,6).Value = "(ATH)"
$ws.Cells.Item(301,7).Value = "Ryanair "
$ws.Cells.Item(301,8).Value = "B738"
$ws.Cells.Item(301,9).Value = "(SP-RSB)"
$ws.Cells.Item(301,10).Value = "12:12 PM"
$ws.Cells.Item(301,12).Value = "0 hours, 7 minutes"

# Row 302 (NUMBER 301)
$ws.Cells.Item(302,1).Value = 301
$ws.Cells.Item(302,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(302,3).Value = "12:20 PM"
$ws.Cells.Item(302,4).Value = "FR7100"
$ws.Cells.Item(302,5).Value = "Oslo"
$ws.Cells.Item(302,6).Value = "(OSL)"
$ws.Cells.Item(302,7).Value = "Ryanair "
$ws.Cells.Item(302,8).Value = "B738"
$ws.Cells.Item(302,9).Value = "(SP-RSH)"
$ws.Cells.Item(302,10).Value = "12:22 PM"
$ws.Cells.Item(302,12).Value = "0 hours, 2 minutes"

# Row 303 (NUMBER 302)
$ws.Cells.Item(303,1).Value = 302
$ws.Cells.Item(303,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(303,3).Value = "12:35 PM"
$ws.Cells.Item(303,4).Value = "W61081"
$ws.Cells.Item(303,5).Value = "Tenerife"
$ws.Cells.Item(303,6).Value = "(TFS)"
$ws.Cells.Item(303,7).Value = "Wizz Air "
$ws.Cells.Item(303,8).Value = "A21N"
$ws.Cells.Item(303,9).Value = "(HA-LZF)"
$ws.Cells.Item(303,10).Value = "12:57 PM"
$ws.Cells.Item(303,12).Value = "0 hours, 22 minutes"

# Row 304 (NUMBER 303)
$ws.Cells.Item(304,1).Value = 303
$ws.Cells.Item(304,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(304,3).Value = "1:40 PM"
$ws.Cells.Item(304,4).Value = "W61049"
$ws.Cells.Item(304,5).Value = "Naples"
$ws.Cells.Item(304,6).Value = "(NAP)"
$ws.Cells.Item(304,7).Value = "Wizz Air "
$ws.Cells.Item(304,8).Value = "A21N"
$ws.Cells.Item(304,9).Value = "(HA-LZJ)"
$ws.Cells.Item(304,10).Value = "1:46 PM"
$ws.Cells.Item(304,12).Value = "0 hours, 6 minutes"

# Row 305 (NUMBER 304)
$ws.Cells.Item(305,1).Value = 304
$ws.Cells.Item(305,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(305,3).Value = "1:45 PM"
$ws.Cells.Item(305,4).Value = "FR5398"
$ws.Cells.Item(305,5).Value = "Dublin"
$ws.Cells.Item(305,6).Value = "(DUB)"
$ws.Cells.Item(305,7).Value = "Ryanair "
$ws.Cells.Item(305,8).Value = "B738"
$ws.Cells.Item(305,9).Value = "(EI-EVC)"
$ws.Cells.Item(305,10).Value = "1:53 PM"
$ws.Cells.Item(305,12).Value = "0 hours, 8 minutes"

# Row 306 (NUMBER 305)
$ws.Cells.Item(306,1).Value = 305
$ws.Cells.Item(306,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(306,3).Value = "2:15 PM"
$ws.Cells.Item(306,4).Value = "LH1357"
$ws.Cells.Item(306,5).Value = "Frankfurt"
$ws.Cells.Item(306,6).Value = "(FRA)"
$ws.Cells.Item(306,7).Value = "Lufthansa "
$ws.Cells.Item(306,8).Value = "CRJ9"
$ws.Cells.Item(306,9).Value = "(D-ACNB)"
$ws.Cells.Item(306,10).Value = "2:20 PM"
$ws.Cells.Item(306,12).Value = "0 hours, 5 minutes"

# Row 307 (NUMBER 306)
$ws.Cells.Item(307,1).Value = 306
$ws.Cells.Item(307,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(307,3).Value = "2:40 PM"
$ws.Cells.Item(307,4).Value = "UNKNOWN"
$ws.Cells.Item(307,5).Value = "Warsaw"
$ws.Cells.Item(307,6).Value = "(WAW)"
$ws.Cells.Item(307,7).Value = "Enter Air "
$ws.Cells.Item(307,8).Value = "B738"
$ws.Cells.Item(307,9).Value = "(SP-ENW)"
$ws.Cells.Item(307,10).Value = "3:59 PM"
$ws.Cells.Item(307,12).Value = "1 hours, 19 minutes"

# Row 308 (NUMBER 307)
$ws.Cells.Item(308,1).Value = 307
$ws.Cells.Item(308,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(308,3).Value = "2:55 PM"
$ws.Cells.Item(308,4).Value = "W61251"
$ws.Cells.Item(308,5).Value = "Larnaca"
$ws.Cells.Item(308,6).Value = "(LCA)"
$ws.Cells.Item(308,7).Value = "Wizz Air "
$ws.Cells.Item(308,8).Value = "A21N"
$ws.Cells.Item(308,9).Value = "(HA-LVT)"
$ws.Cells.Item(308,10).Value = "2:56 PM"
$ws.Cells.Item(308,12).Value = "0 hours, 1 minutes"

# Row 309 (NUMBER 308)
$ws.Cells.Item(309,1).Value = 308
$ws.Cells.Item(309,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(309,3).Value = "3:10 PM"
$ws.Cells.Item(309,4).Value = "LO3884"
$ws.Cells.Item(309,5).Value = "Warsaw"
$ws.Cells.Item(309,6).Value = "(WAW)"
$ws.Cells.Item(309,7).Value = "LOT "
$ws.Cells.Item(309,8).Value = "E75S"
$ws.Cells.Item(309,9).Value = "(SP-LID)"
$ws.Cells.Item(309,10).Value = "3:10 PM"
$ws.Cells.Item(309,12).Value = "0 hours, 0 minutes"

Write-Output "done"